$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.011370182037354
$ws.Range("B1").Value = 2.195341825485229
$ws.Range("C1").Value = 2.49753212928772
$ws.Range("D1").Value = 3.804660081863403
$ws.Range("E1").Value = 1.260306239128113
